$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "FuncLoc" and "Previous Doc" validation data values on row 2
$ws.Range("AV2").Value = "ASMPD503ALVINCOMMUN"
$ws.Range("AX2").Value = "2152430001"
